$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-30 all get Area -> "SLO" and Station -> "45bis"
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 2).Value = "SLO"
    $ws.Cells.Item($r, 3).Value = "45bis"
}
